{"js": "// Update the date heading and the worksheet's answer cells.\n//\n// The document is a single title paragraph (\"YYYY-MM-DD Weekday\") followed\n// by one table whose rows alternate between a populated \"answers\" row\n// (5 cells of \"A\u00f7B=Q, R\" text) and 3 blank spacer rows. Every populated\n// cell's text changes value (but no rows/cells/paragraphs are added or\n// removed), so the edit is a same-shape text-content replacement driven by\n// positional (row-major, left-to-right) order.\n\n// 1. Update the date/weekday title (first paragraph of the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titleParagraph = paragraphs.items[0];\ntitleParagraph.insertText(\"2024-03-21 Thursday\", \"Replace\");\n\n// 2. Update the table's answer cells, preserving every row's formatting\n//    (table.values round-trips cell text only, leaving run/paragraph\n//    properties such as fonts untouched).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// New answer strings, in row-major reading order, for every non-blank row.\nconst newAnswers = [\n  [\"55\u00f73=18, 1\", \"89\u00f77=12, 5\", \"83\u00f79=9, 2\", \"46\u00f76=7, 4\", \"61\u00f78=7, 5\"],\n  [\"84\u00f74=21, 0\", \"13\u00f79=1, 4\", \"44\u00f77=6, 2\", \"51\u00f72=25, 1\", \"82\u00f77=11, 5\"],\n  [\"59\u00f79=6, 5\", \"80\u00f72=40, 0\", \"41\u00f75=8, 1\", \"77\u00f73=25, 2\", \"34\u00f75=6, 4\"],\n  [\"69\u00f78=8, 5\", \"44\u00f78=5, 4\", \"86\u00f76=14, 2\", \"55\u00f75=11, 0\", \"60\u00f72=30, 0\"],\n  [\"46\u00f76=7, 4\", \"27\u00f75=5, 2\", \"44\u00f78=5, 4\", \"74\u00f75=14, 4\", \"89\u00f74=22, 1\"],\n];\n\nconst updatedValues = table.values.map((row) => row.slice());\nlet answerRowIdx = 0;\nfor (let r = 0; r < updatedValues.length; r++) {\n  const isAnswerRow = updatedValues[r].some((cell) => cell !== \"\");\n  if (isAnswerRow) {\n    updatedValues[r] = newAnswers[answerRowIdx].slice();\n    answerRowIdx++;\n  }\n}\n\ntable.values = updatedValues;\nawait context.sync();\n", "ps1": "# Update the date heading and the worksheet's answer cells.\n#\n# The document is a single title paragraph (\"YYYY-MM-DD Weekday\") followed\n# by one table whose rows alternate between a populated \"answers\" row\n# (5 cells of \"A\u00f7B=Q, R\" text) and 3 blank spacer rows. Every populated\n# cell's text changes value (but no rows/cells/paragraphs are added or\n# removed), so the edit is a same-shape text-content replacement driven by\n# positional (row-major, left-to-right) order.\n\n$d = $word.ActiveDocument\n\n# 1. Update the date/weekday title (first paragraph of the body).\n$d.Paragraphs.Item(1).Range.Text = \"2024-03-21 Thursday\"\n\n# 2. Update the table's answer cells, preserving every row's formatting\n#    (assigning Cell.Range.Text only swaps the run's text, leaving\n#    paragraph/run properties such as fonts untouched).\n$t = $d.Tables.Item(1)\n\n# New answer strings, in row-major reading order, for every non-blank row.\n$newAnswers = @(\n    @(\"55\u00f73=18, 1\", \"89\u00f77=12, 5\", \"83\u00f79=9, 2\", \"46\u00f76=7, 4\", \"61\u00f78=7, 5\"),\n    @(\"84\u00f74=21, 0\", \"13\u00f79=1, 4\", \"44\u00f77=6, 2\", \"51\u00f72=25, 1\", \"82\u00f77=11, 5\"),\n    @(\"59\u00f79=6, 5\", \"80\u00f72=40, 0\", \"41\u00f75=8, 1\", \"77\u00f73=25, 2\", \"34\u00f75=6, 4\"),\n    @(\"69\u00f78=8, 5\", \"44\u00f78=5, 4\", \"86\u00f76=14, 2\", \"55\u00f75=11, 0\", \"60\u00f72=30, 0\"),\n    @(\"46\u00f76=7, 4\", \"27\u00f75=5, 2\", \"44\u00f78=5, 4\", \"74\u00f75=14, 4\", \"89\u00f74=22, 1\")\n)\n\n$answerRowIdx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $raw = $t.Cell($r, 1).Range.Text\n    $core = $raw.Substring(0, $raw.Length - 2)   # strip trailing cell-mark (\\r\\a)\n    if ($core.Length -gt 0) {\n        $rowValues = $newAnswers[$answerRowIdx]\n        for ($c = 1; $c -le $t.Columns.Count; $c++) {\n            $t.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n        }\n        $answerRowIdx++\n    }\n}\n"}
